# Update the "Duplicate_Management" sheet:
#  - B2 gets a new URL and the "Hyperlink" visual style (underline + theme
#    hyperlink colour) without an actual live hyperlink
#  - a new column C ("Action Value") is added with a header + a text "1"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Duplicate_Management")
$ws.Activate()

# --- B2: replace the URL text and give it the Hyperlink look -------------
$ws.Range("B2").Value = "https://www.google.com/maps/place/?q=place_id:ChIJLdD4Ikv1xokRXCXMOatucy4"

# Adding then deleting a real hyperlink leaves the cell formatted with the
# built-in "Hyperlink" cell style (underline, theme hyperlink colour) while
# removing the live hyperlink / relationship again, matching the workbook.
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.google.com/maps/place/?q=place_id:ChIJLdD4Ikv1xokRXCXMOatucy4") | Out-Null
$ws.Hyperlinks.Delete()

# --- New column C: "Action Value" header + value --------------------------
$ws.Range("C1").Value = "Action Value"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> bold header style
$excel.CutCopyMode = $false

$ws.Columns.Item(3).ColumnWidth = 11.5

# Store the value as text "1" (quote-prefixed), matching the source data
$ws.Range("C2").Value = "'1"

# --- Selection housekeeping ------------------------------------------------
$ws.Range("C3").Select() | Out-Null
